# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value, identical for both sheets except rows 20-25
# shift by +1 on "全部类型" because it contains two extra rows.
$updatesExhibition = @{
    2  = 7061
    4  = 75
    7  = 7560
    8  = 81
    13 = 426
    14 = 158
    16 = 423
    19 = 22
    20 = 5399
    21 = 137
    22 = 192
    23 = 833
    25 = 281
}

$updatesAllTypes = @{
    2  = 7061
    4  = 75
    7  = 7560
    8  = 81
    13 = 426
    14 = 158
    16 = 423
    19 = 22
    21 = 5399
    23 = 137
    24 = 192
    25 = 833
    27 = 281
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $updatesAllTypes[$row]
}
